$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.726.53"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "2.523.67"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.10%  "
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.532"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.77"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0810"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.56"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.31%  "
$ws.Range("E13").Value = "  -2.66%  "
$ws.Range("D14").Value = "2.914.29"
$ws.Range("E14").Value = "  +0.47%  "
$ws.Range("D15").Value = "2.526.85"
$ws.Range("E15").Value = "  -0.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.846"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.13%  "
$ws.Range("D18").Value = "42.809.88"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.96%  "
$ws.Range("D21").Value = "0.0₃0960"
$ws.Range("E21").Value = "  -0.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "251.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.33%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("E28").Value = "  +2.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.90"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.22%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.72"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.91%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.90"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("E33").Value = "  +3.77%  "
$ws.Range("E34").Value = "  +4.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.31"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.82"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0789"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.71%  "
$ws.Range("E38").Value = "  -0.18%  "
$ws.Range("B39").Value = "ApeXProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.41"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +14.97%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.118"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.94%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "21.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -9.43%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0305"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.67%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.81"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.66%  "
$ws.Range("E44").Value = "  +0.29%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.024.39"
$ws.Range("E45").Value = "  -2.23%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.25"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.05"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "84.02"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.49%  "
$ws.Range("D51").Value = "2.768.82"
$ws.Range("E51").Value = "  +0.34%  "
